$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Composite + Original Feedstocks"

$data = @(
  @("DCW-BSG, 190°C, 1hr", 0.000084485016096908977282185782, 0.030436843757038738872422456438, 0.932647287161454396908766284469),
  @("SRU-DCW, 190°C, 1hr", 0.000103199446193244203100365664, 0.021928682096168050641438185266, 0.932279257564905039501468309027),
  @("SRU-DCW-BSG, 190°C, 1hr", 0.000104351348826945699329975248, 0.021541143842139020420090034236, 0.931054199447511443565872468753),
  @("BSG, 190°C, 1hr", 0.000106339630501552402240712791, 0.030186207721261491127018317115, 0.909812848633535176112729914166),
  @("SRU-BSG, 190°C, 1hr", 0.000123459600539317599981678497, 0.018947057960492081490899707319, 0.904661498251477125265296308498),
  @("DCW-BSG, 220°C, 1hr", 0.000099630519818505412494961948, 0.034171579072517471731273985824, 0.904443008510913193376268282009),
  @("SRU-DCW, 220°C, 1hr", 0.0001233051920642406096718241, 0.024472485695532029958965125616, 0.899275486669544266149500799656),
  @("SRU-DCW-BSG, 220°C, 1hr", 0.00012360817168698460220832569, 0.024214798408418409386699821084, 0.899223867293164946801198311732),
  @("DCW-BSG, 250°C, 1hr", 0.000115341251308633305867543117, 0.037749426885934467268413072816, 0.875540522393431941416963582014),
  @("SRU, 190°C, 1hr", 0.000143786564268074887329440514, 0.018147346553377241562543531472, 0.875517375004525866266646971781),
  @("SRU-DCW-BSG, 250°C, 1hr", 0.000143866514012248490009954205, 0.026815615024599240290292812006, 0.866321927517144318997566188045),
  @("SRU-DCW, 250°C, 1hr", 0.000144640024429078495023584527, 0.02693027958529788101760260588, 0.865048491738738101197725427483),
  @("SRU-BSG, 220°C, 1hr", 0.000151132365637423606800210374, 0.020976577250358049986367348083, 0.863186267615368429950706286036),
  @("BSG, 220°C, 1hr", 0.000141735199977738698029136599, 0.030538787105512308489574380133, 0.862329630487106757463777739758),
  @("BSG, 250°C, 1hr", 0.000154917280283288001001620882, 0.034282151512224696976449678232, 0.836423667475505117074874306127),
  @("SRU-BSG, 250°C, 1hr", 0.000181223516128760603165204213, 0.022958283078932439913444341073, 0.819434895497955118059962842381),
  @("SRU, 220°C, 1hr", 0.000187880932093521612171252433, 0.020141919103052168676937228042, 0.813165290630568193819271982647),
  @("DCW-BSG, 190°C, 3hr", 0.000168624886278890911105177164, 0.043425935403065440976178024357, 0.795490581766866267088289532694),
  @("SRU, 250°C, 1hr", 0.000211299698264508792680538285, 0.022178560815777209613841591818, 0.78044938414750120170282343679),
  @("DCW-BSG, 220°C, 3hr", 0.000157616853562017497332806149, 0.056431968782394502370802769065, 0.766524288458935609824607126939),
  @("SRU-DCW, 190°C, 3hr", 0.000235206618322910001152650783, 0.037726797994676770919841857221, 0.725638188943979245415505374694),
  @("SRU-DCW-BSG, 190°C, 3hr", 0.000237743640814324893109632519, 0.037494966810220686781462973158, 0.722990833988270131804654283769),
  @("SRU-DCW, 220°C, 3hr", 0.000228199613781797903695053309, 0.0460011194623969194017476525, 0.715556701228577729523294692626),
  @("SRU-DCW-BSG, 220°C, 3hr", 0.000229522666124418092117517864, 0.045732447951780212658245972079, 0.714618163094369762688984337728),
  @("DCW-BSG, 250°C, 3hr", 0.000196758030530606989111017113, 0.06123958381921436322503637939, 0.708163504760110362745706424903),
  @("BSG, 190°C, 3hr", 0.000233461052364261303264728054, 0.050163998139006522880567473521, 0.698611306313244306842591413442),
  @("SRU-BSG, 190°C, 3hr", 0.000281451090838755299712581337, 0.035904033043353607290715245881, 0.6745314009936936860967193752),
  @("BSG, 220°C, 3hr", 0.000263963824077419076472722503, 0.055782180695747148047036745311, 0.648264764342243449135594346444),
  @("SRU-DCW-BSG, 250°C, 3hr", 0.000285339183604000994325805562, 0.050878052223307532186957757858, 0.637234181140629130091213028209),
  @("SRU-DCW, 250°C, 3hr", 0.000285686456013399894555276903, 0.051045193738535367400377396052, 0.636434620369530024142079582816),
  @("SRU-BSG, 220°C, 3hr", 0.000306236140810142178601827023, 0.042207916440134829461072740742, 0.634662486629954281092125256691),
  @("SRU, 190°C, 3hr", 0.000319943916139601026559641239, 0.036281510968898818270123030061, 0.632346460248279007565486153908),
  @("BSG, 250°C, 3hr", 0.000296620816986829522864804742, 0.059266601688229857136658296213, 0.602982831321567047666576399934),
  @("SRU, 220°C, 3hr", 0.000400690257657189408124764896, 0.040872858835960167167566936541, 0.549769136479078213142202002928),
  @("SRU-BSG, 250°C, 3hr", 0.000387973102116955686265858949, 0.047149173561838682300351166532, 0.546431875219450802205756190233),
  @("DCW, 220°C, 1hr", 0.000076755697230615900191266765, 0.147727666590357187281767892273, 0.516597225943495441669028878096),
  @("DCW, 190°C, 3hr", 0.00007705723325053621323814379, 0.149735065942135592287343115458, 0.511615510142476859201110528375),
  @("SRU, 250°C, 3hr", 0.00046234366521296109192248891, 0.045183972918843316735060255951, 0.497841742847301815810823200081),
  @("DCW, 250°C, 1hr", 0.000088852952615613801706437314, 0.164027851961231790900441751546, 0.472530618429845394867783170412),
  @("DCW, 190°C, 1hr", 0.000061521644857058237101758713, 0.179218741620582100448899609546, 0.462892988178742914584518075571),
  @("DCW, 220°C, 3hr", 0.000089333895927881830367756066, 0.184427370489045394652549703096, 0.436583411281966204953164378821),
  @("DCW, 250°C, 3hr", 0.00010600369791037919523293831, 0.185264129887337591906515399387, 0.423517299266585878214641525119)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
